# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $TextValue) {
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $TextValue
    $r.Style = "Normal"
}

Set-TextValue 'D2' '67.340.74'
Set-TextValue 'E2' '  +0.58%  '
Set-TextValue 'D3' '3.113.40'
Set-TextValue 'E3' '  +1.19%  '
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '575.64'
Set-TextValue 'E5' '  -0.28%  '
Set-TextValue 'D6' '178.05'
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '3.114.18'
Set-TextValue 'E8' '  +1.32%  '
Set-TextValue 'D9' '0.516'
Set-TextValue 'E9' '  +0.62%  '
Set-TextValue 'D10' '6.53'
Set-TextValue 'E10' '  +2.44%  '
Set-TextValue 'E11' '  +1.25%  '
Set-TextValue 'D12' '0.468'
Set-TextValue 'E12' '  -0.63%  '
Set-TextValue 'D13' '0.0000242'
Set-TextValue 'E13' '  +0.13%  '
Set-TextValue 'D14' '36.49'
Set-TextValue 'E14' '  +1.39%  '
Set-TextValue 'E15' '  +0.89%  '
Set-TextValue 'D16' '3.628.49'
Set-TextValue 'E16' '  +1.05%  '
Set-TextValue 'D17' '67.270.45'
Set-TextValue 'E17' '  +0.61%  '
Set-TextValue 'D18' '7.04'
Set-TextValue 'E18' '  +0.17%  '
Set-TextValue 'D19' '3.111.65'
Set-TextValue 'E19' '  +1.20%  '
Set-TextValue 'D20' '16.51'
Set-TextValue 'E20' '  -2.54%  '
Set-TextValue 'D21' '486.92'
Set-TextValue 'E21' '  +0.06%  '
Set-TextValue 'B22' 'Polygon'
Set-TextValue 'C22' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D22' '0.689'
Set-TextValue 'E22' '  -0.10%  '
Set-TextValue 'B23' 'Uniswap'
Set-TextValue 'C23' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D23' '7.71'
Set-TextValue 'E23' '  +0.12%  '
Set-TextValue 'D24' '83.68'
Set-TextValue 'E24' '  +1.13%  '
Set-TextValue 'D25' '12.79'
Set-TextValue 'E25' '  -0.29%  '
Set-TextValue 'D26' '2.27'
Set-TextValue 'E26' '  +2.40%  '
Set-TextValue 'D27' '10.46'
Set-TextValue 'E27' '  +1.28%  '
Set-TextValue 'E28' '  +0.09%  '
Set-TextValue 'D29' '7.95'
Set-TextValue 'E29' '  +1.65%  '
Set-TextValue 'E30' '  +1.28%  '
Set-TextValue 'E31' '  -0.27%  '
Set-TextValue 'D32' '28.13'
Set-TextValue 'E32' '  +1.71%  '
Set-TextValue 'E33' '  +0.59%  '
Set-TextValue 'E34' '  +3.25%  '
Set-TextValue 'E35' '  +0.01%  '
Set-TextValue 'D36' '47.94'
Set-TextValue 'E36' '  +3.55%  '
Set-TextValue 'D37' '0.949'
Set-TextValue 'E37' '  -0.28%  '
Set-TextValue 'D38' '5.59'
Set-TextValue 'E38' '  -1.17%  '
Set-TextValue 'D39' '0.318'
Set-TextValue 'E39' '  +5.30%  '
Set-TextValue 'D40' '49.24'
Set-TextValue 'D41' '2.01'
Set-TextValue 'E41' '  +1.27%  '
Set-TextValue 'E42' '  +0.28%  '
Set-TextValue 'D43' '8.30'
Set-TextValue 'E43' '  -0.40%  '
Set-TextValue 'D45' '2.787.58'
Set-TextValue 'E45' '  +1.03%  '
Set-TextValue 'D46' '373.64'
Set-TextValue 'E46' '  +0.69%  '
Set-TextValue 'B47' 'InjectiveProtocol'
Set-TextValue 'C47' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D47' '26.75'
Set-TextValue 'E47' '  +9.19%  '
Set-TextValue 'B48' 'VeChain'
Set-TextValue 'C48' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D48' '0.0346'
Set-TextValue 'E48' '  +0.33%  '
Set-TextValue 'B49' 'Monero'
Set-TextValue 'C49' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D49' '135.76'
Set-TextValue 'E49' '  -0.30%  '
Set-TextValue 'D51' '2.35'
Set-TextValue 'E51' '  +9.11%  '

Write-Host "Applied 93 cell updates to Sheet1"
